# Applies the "Client Information" content edits described by the commit:
#   - dictionary cleanup: drop HOME TEL / BUSSINESS TEL rows, rename "MOBILE "
#     -> "MOBILE", reorder the phone block on the Info sheet
#   - split the combined residential address into LINE 1 / LINE 2
#   - fix "Maried " -> "Married"
#   - replace the placeholder beneficiary rows (a/b/c/d/e, k/j/h/g/f) with a
#     real beneficiary (TOLANI SANJAY, spouse, 100%) on the Beneficiaries sheet

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Info" sheet
# ---------------------------------------------------------------------------
$info = $wb.Worksheets.Item("Info")

# Drop the "HOME TEL" (row 21) and "BUSSINESS TEL" (row 22) rows entirely;
# everything below shifts up by two rows.
$info.Rows.Item(21).Delete()
$info.Rows.Item(21).Delete()

# The phone/contact block (now rows 20-25) is reordered to:
#   COUNTRY CODE, AREA CODE, RESIDENTIAL NUMBER, MOBILE, OFFICE NUMBER, EMAIL
$info.Cells.Item(20, 1).Value = "COUNTRY CODE"
$info.Cells.Item(20, 2).Value = "971"

$info.Cells.Item(21, 1).Value = "AREA CODE"
$info.Cells.Item(21, 2).Value = "0"

$info.Cells.Item(22, 1).Value = "RESIDENTIAL NUMBER"
$info.Cells.Item(22, 2).Value = "0"

$info.Cells.Item(23, 1).Value = "MOBILE"
$info.Cells.Item(23, 2).Value = "553595566"

$info.Cells.Item(24, 1).Value = "OFFICE NUMBER"
$info.Cells.Item(24, 2).Value = "0"

$info.Cells.Item(25, 1).Value = "EMAIL "
$info.Cells.Item(25, 2).Value = "gib@eim.ae"

# Split the previously-combined residential address into its two lines.
$info.Cells.Item(31, 2).Value = "Villa-18/2A, 394/Emirates Hill, Third (Meadows-8), "
$info.Cells.Item(32, 2).Value = "Premise Number 394041593, Dubai, UAE"

# Typo fix.
$info.Cells.Item(40, 2).Value = "Married"

# ---------------------------------------------------------------------------
# "Beneficiaries" sheet
# ---------------------------------------------------------------------------
$ben = $wb.Worksheets.Item("Beneficiaries")

$ben.Cells.Item(3, 2).Value = "TOLANI SANJAY"
$ben.Cells.Item(3, 3).Value = "(SPOUSE)"
$ben.Cells.Item(3, 4).Value = "(SPOUSE)"
$ben.Cells.Item(3, 5).Value = ""
$ben.Cells.Item(3, 6).Value = "100%"

# Remove the old placeholder dictionary rows (a/b/c/d/e and k/j/h/g/f).
$ben.Rows.Item(4).ClearContents()
$ben.Rows.Item(5).ClearContents()
